$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.966.97'
$ws.Range("E2").Value = '  +2.10%  '

$ws.Range("D3").Value = '1.700.56'
$ws.Range("E3").Value = '  +0.81%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").Value = '''315.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.20%  '

$ws.Range("E6").Value = '  +0.28%  '

$ws.Range("D7").Value = '''0.3978'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.70%  '

$ws.Range("D8").Value = '''0.4025'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("E9").Value = '  -1.60%  '

$ws.Range("D10").Value = '''52.93'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.99%  '

$ws.Range("D11").Value = '''1.002'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.39%  '

$ws.Range("D12").Value = '''0.08797'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.30%  '

$ws.Range("D13").Value = '''26.04'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.76%  '

$ws.Range("D14").Value = '''7.465'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.11%  '

$ws.Range("D15").Value = '''0.00001350'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.68%  '

$ws.Range("D16").Value = '''7.945'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.80%  '

$ws.Range("D17").Value = '1.708.45'
$ws.Range("E17").Value = '  +2.27%  '

$ws.Range("D18").Value = '''96.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.02%  '

$ws.Range("D19").Value = '''0.07199'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.48%  '

$ws.Range("D20").Value = '''20.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.55%  '

$ws.Range("D21").Value = '''7.320'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.69%  '

$ws.Range("D22").Value = '''1.002'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.21%  '

$ws.Range("D23").Value = '''14.38'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.73%  '

$ws.Range("D24").Value = '24.968.72'
$ws.Range("E24").Value = '  +2.10%  '

$ws.Range("D25").Value = '''2.357'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.00%  '

$ws.Range("D26").Value = '''2.940'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.13%  '

$ws.Range("D27").Value = '''23.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.24%  '

$ws.Range("D28").Value = '''6.199'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +16.14%  '

$ws.Range("D29").Value = '''162.06'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.19%  '

$ws.Range("D30").Value = '''151.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.61%  '

$ws.Range("D31").Value = '''8.349'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.52%  '

$ws.Range("D32").Value = '''2.631'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +25.66%  '

$ws.Range("D33").Value = '1.896.85'
$ws.Range("E33").Value = '  +2.01%  '

$ws.Range("D34").Value = '''0.08576'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.22%  '

$ws.Range("D35").Value = '''0.03139'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.51%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''1.039'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.68%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '''7.156'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.69%  '

$ws.Range("D38").Value = '''0.2877'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.23%  '

$ws.Range("D39").Value = '''0.09575'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.81%  '

$ws.Range("E40").Value = '  +0.48%  '

$ws.Range("D41").Value = '''0.8250'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.83%  '

$ws.Range("D42").Value = '''14.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.04%  '

$ws.Range("D43").Value = '''1.485'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.72%  '

$ws.Range("E44").Value = '  -1.70%  '

$ws.Range("D45").Value = '''2.681'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.90%  '

$ws.Range("D46").Value = '''0.7386'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.04%  '

$ws.Range("D47").Value = '''4.242'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.52%  '

$ws.Range("E48").Value = '  -1.23%  '

$ws.Range("D49").Value = '''0.08785'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.71%  '

$ws.Range("D50").Value = '''1.002'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.21%  '

$ws.Range("D51").Value = '''139.23'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.03%  '
